$d = $word.ActiveDocument

# The two "Intermediate Precision" style tables (Sample / n / Mean / SD / CV)
# that need their literal values swapped for Jinja template placeholders,
# and their centered paragraph justification stripped (both on the index
# column and the data columns).
$targetTables = @(5, 6)

# Column index (2-5) -> placeholder field-name suffix. Column 1 holds the
# plain "1"/"2"/"3" sample index and keeps its text (only loses centering).
$fieldNames = @("n", "mean", "sd", "cv")

foreach ($tblIndex in $targetTables) {
    $tbl = $d.Tables.Item($tblIndex)

    for ($sampleNum = 1; $sampleNum -le 3; $sampleNum++) {
        $rowIndex = $sampleNum + 1   # row 1 is the header row

        # Column 1: the plain "1"/"2"/"3" sample-index cell - just drop the
        # centered justification, leave the text alone.
        $idxCell = $tbl.Cell($rowIndex, 1)
        $idxCell.Range.Paragraphs.Item(1).Alignment = 0

        # Columns 2-5: n / mean / sd / cv - replace the literal value with
        # the templated placeholder and drop the centered justification.
        for ($col = 2; $col -le 5; $col++) {
            $fieldName = $fieldNames[$col - 2]
            $placeholder = "{{ inter_var_sample" + $sampleNum + "_" + $fieldName + " }}"

            $dataCell = $tbl.Cell($rowIndex, $col)
            $textRange = $dataCell.Range
            # Cell.Range.Text includes a trailing paragraph mark + cell mark
            # (chr(13) + chr(7)); shrink the range so we only overwrite the
            # actual visible run text, not the table-cell delimiters.
            $textRange.End = $textRange.End - 2
            $textRange.Text = $placeholder

            $dataCell.Range.Paragraphs.Item(1).Alignment = 0
        }
    }
}
